# Auto-ish generated script applying the FujitsuTiers.xlsx edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural: insert the two blank separator rows ---
# (pushes Group-2 block down by one row, then Group-3 block down by one more)
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(24).Insert()

# --- Header row: new columns C..F ---
$c = $ws.Range("C1")
$c.NumberFormat = "@"
$c.HorizontalAlignment = -4131
$c.Value = "Before dir"
$c = $ws.Range("D1")
$c.NumberFormat = "@"
$c.HorizontalAlignment = -4131
$c.Value = "After dir"

$ws.Range("E1").Value = "Before Notes"
$ws.Range("F1").Value = "After notes"
$ws.Range("E2").Value = "2nd epoch is bad here"
$ws.Range("E3").Value = "OK"
$ws.Range("E4").Value = "OK"
$ws.Range("E5").Value = "bad channels here"
$ws.Range("E6").Value = "1st epoch is bad here"

# --- "After dir"/"Before dir" data cells: text-typed (kept as text, format "@") ---
$c = $ws.Range("C2")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("C3")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("C4")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("C5")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11052018"
$c = $ws.Range("C6")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "09062018"
$c = $ws.Range("C8")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11142018"
$c = $ws.Range("C9")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("C10")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11122018"
$c = $ws.Range("C11")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11012018"
$c = $ws.Range("C12")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("C14")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11132018"
$c = $ws.Range("C15")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11022018"
$c = $ws.Range("C16")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11122018"
$c = $ws.Range("C17")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11092018"
$c = $ws.Range("C18")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11012018"
$c = $ws.Range("C19")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11132018"
$c = $ws.Range("C20")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11092018"
$c = $ws.Range("C21")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11072018"
$c = $ws.Range("C22")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11132018"
$c = $ws.Range("C23")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11072018"
$c = $ws.Range("C24")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("C25")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("C26")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152

# --- "After dir"/"Before dir" data cells: number-typed (format General) ---
$c = $ws.Range("D2")
$c.NumberFormat = "General"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11062018"
$c = $ws.Range("D3")
$c.NumberFormat = "General"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11052018"
$c = $ws.Range("D4")
$c.NumberFormat = "General"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11062018"
$c = $ws.Range("D6")
$c.NumberFormat = "General"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11062018"
$c = $ws.Range("C7")
$c.NumberFormat = "General"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "9072018"
$c = $ws.Range("D7")
$c.NumberFormat = "General"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11072018"
$c = $ws.Range("D9")
$c.NumberFormat = "General"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11072018"
$c = $ws.Range("C13")
$c.NumberFormat = "General"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "9052018"
$c = $ws.Range("D13")
$c.NumberFormat = "General"
$c.Font.Size = 9
$c.HorizontalAlignment = -4152
$c.Value = "11142018"

